$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1216.7778
$ws.Range("I43").Value = 1050
$ws.Range("J43").Value = 1350.2
$ws.Range("K43").Value = 1050
$ws.Range("L43").Value = 1350.2
$ws.Range("M43").Value = -981
$ws.Range("N43").Value = -1488.2

$ws.Range("H76").Value = 3499.5
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3499.5
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 3499.5
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -4129.5

$ws.Range("H79").Value = 3499.5
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3499.5
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 3499.5
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -5683.5

$ws.Range("H132").Value = 1211.5156
$ws.Range("I132").Value = 1171.8545
$ws.Range("J132").Value = 1453.8889
$ws.Range("K132").Value = 3515.5635
$ws.Range("L132").Value = 4361.6667
$ws.Range("M132").Value = -985.5634999999997
$ws.Range("N132").Value = -9421.6667

$ws.Range("H137").Value = 1121.317
$ws.Range("I137").Value = 938.28
$ws.Range("J137").Value = 1407.3125
$ws.Range("K137").Value = 2814.84
$ws.Range("L137").Value = 4221.9375
$ws.Range("M137").Value = -264.8400000000001
$ws.Range("N137").Value = -9321.9375

$ws.Range("H141").Value = 876917.5
$ws.Range("I141").Value = 1037914.56
$ws.Range("J141").Value = 7533.2
$ws.Range("K141").Value = 3113743.68
$ws.Range("L141").Value = 22599.6
$ws.Range("M141").Value = -3108563.68
$ws.Range("N141").Value = -32959.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3158.1667
$ws.Range("I32").Value = 2678.7
$ws.Range("J32").Value = 12747.5
$ws.Range("K32").Value = 2678.7
$ws.Range("L32").Value = 12747.5
$ws.Range("M32").Value = -2391.7
$ws.Range("N32").Value = -13321.5

$ws.Range("H61").Value = 1969.0541
$ws.Range("I61").Value = 1453
$ws.Range("J61").Value = 11000
$ws.Range("K61").Value = 1453
$ws.Range("L61").Value = 11000
$ws.Range("M61").Value = -1241
$ws.Range("N61").Value = -11424

$ws.Range("H74").Value = 1254.439
$ws.Range("I74").Value = 964.5
$ws.Range("K74").Value = 964.5
$ws.Range("M74").Value = -90.5

$ws.Range("H77").Value = 1254.439
$ws.Range("I77").Value = 964.5
$ws.Range("K77").Value = 4822.5
$ws.Range("M77").Value = -454.5

$ws.Range("H122").Value = 1573.6897
$ws.Range("J122").Value = 1893.375
$ws.Range("L122").Value = 5680.125
$ws.Range("N122").Value = -10580.125

$ws.Range("H125").Value = 50500
$ws.Range("J125").Value = 50500
$ws.Range("L125").Value = 50500
$ws.Range("N125").Value = -60340

$ws.Range("H132").Value = 1199.5802
$ws.Range("I132").Value = 894.6229
$ws.Range("K132").Value = 2683.8687
$ws.Range("M132").Value = -153.8687

$ws.Range("H136").Value = 1969.0541
$ws.Range("I136").Value = 1453
$ws.Range("J136").Value = 11000
$ws.Range("K136").Value = 4359
$ws.Range("L136").Value = 33000
$ws.Range("M136").Value = -1809
$ws.Range("N136").Value = -38100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 488.45456
$ws.Range("I94").Value = 559.625
$ws.Range("K94").Value = 559.625
$ws.Range("M94").Value = -108.625

$ws.Range("H134").Value = 3297.7659
$ws.Range("I134").Value = 3217.1777
$ws.Range("J134").Value = 5111
$ws.Range("K134").Value = 9651.533100000001
$ws.Range("L134").Value = 15333
$ws.Range("M134").Value = -7116.533100000001
$ws.Range("N134").Value = -20403

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1769.4242
$ws.Range("I31").Value = 1378.0435
$ws.Range("J31").Value = 2669.6
$ws.Range("K31").Value = 1378.0435
$ws.Range("L31").Value = 2669.6
$ws.Range("M31").Value = -1083.0435
$ws.Range("N31").Value = -3259.6

$ws.Range("H34").Value = 1769.4242
$ws.Range("I34").Value = 1378.0435
$ws.Range("J34").Value = 2669.6
$ws.Range("K34").Value = 1378.0435
$ws.Range("L34").Value = 2669.6
$ws.Range("M34").Value = -1176.0435
$ws.Range("N34").Value = -3073.6

$ws.Range("H58").Value = 870576.4
$ws.Range("I58").Value = 1208368.4
$ws.Range("J58").Value = 1968.3572
$ws.Range("K58").Value = 1208368.4
$ws.Range("L58").Value = 1968.3572
$ws.Range("M58").Value = -1208165.4
$ws.Range("N58").Value = -2374.3572

$ws.Range("H69").Value = 64095.75
$ws.Range("I69").Value = 12060.667
$ws.Range("K69").Value = 12060.667
$ws.Range("M69").Value = -11311.667

$ws.Range("H72").Value = 64095.75
$ws.Range("I72").Value = 12060.667
$ws.Range("K72").Value = 36182.001
$ws.Range("M72").Value = -32438.001

$ws.Range("H132").Value = 1367.9286
$ws.Range("I132").Value = 1125
$ws.Range("J132").Value = 2171.4614
$ws.Range("K132").Value = 3375
$ws.Range("L132").Value = 6514.3842
$ws.Range("M132").Value = -845
$ws.Range("N132").Value = -11574.3842

$ws.Range("H134").Value = 1531.1608
$ws.Range("I134").Value = 1431.3489
$ws.Range("J134").Value = 1861.3077
$ws.Range("K134").Value = 4294.0467
$ws.Range("L134").Value = 5583.9231
$ws.Range("M134").Value = -1759.0467
$ws.Range("N134").Value = -10653.9231

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 870576.4
$ws.Range("I136").Value = 1208368.4
$ws.Range("J136").Value = 1968.3572
$ws.Range("K136").Value = 3625105.2
$ws.Range("L136").Value = 5905.071599999999
$ws.Range("M136").Value = -3622555.2
$ws.Range("N136").Value = -11005.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10816.829
$ws.Range("I131").Value = 772.25
$ws.Range("J131").Value = 11331.936
$ws.Range("K131").Value = 2316.75
$ws.Range("L131").Value = 33995.808
$ws.Range("M131").Value = 2723.25
$ws.Range("N131").Value = -44075.808

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1298.52
$ws.Range("I122").Value = 1208.7
$ws.Range("K122").Value = 3626.1
$ws.Range("M122").Value = -1176.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2831.2632
$ws.Range("I7").Value = 2089.4
$ws.Range("K7").Value = 2089.4
$ws.Range("M7").Value = -1977.4

$ws.Range("H122").Value = 3060.5715
$ws.Range("I122").Value = 1792.875
$ws.Range("K122").Value = 5378.625
$ws.Range("M122").Value = -2928.625

$ws.Range("H126").Value = 2831.2632
$ws.Range("I126").Value = 2089.4
$ws.Range("K126").Value = 6268.200000000001
$ws.Range("M126").Value = -3798.200000000001

$ws.Range("H132").Value = 1535.9333
$ws.Range("I132").Value = 1235.9445
$ws.Range("K132").Value = 3707.8335
$ws.Range("M132").Value = -1177.8335

$ws.Range("H136").Value = 1977.4058
$ws.Range("I136").Value = 1203.0566
$ws.Range("K136").Value = 3609.1698
$ws.Range("M136").Value = -1059.1698

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 415.25
$ws.Range("I113").Value = 231.55556
$ws.Range("J113").Value = 966.3333
$ws.Range("K113").Value = 694.66668
$ws.Range("L113").Value = 2898.9999
$ws.Range("M113").Value = 1475.33332
$ws.Range("N113").Value = -7238.9999

$ws.Range("H132").Value = 1090.6833
$ws.Range("I132").Value = 824.6818
$ws.Range("J132").Value = 1822.1875
$ws.Range("K132").Value = 2474.0454
$ws.Range("L132").Value = 5466.5625
$ws.Range("M132").Value = 55.95460000000003
$ws.Range("N132").Value = -10526.5625

$ws.Range("H136").Value = 11822516
$ws.Range("I136").Value = 17923740
$ws.Range("J136").Value = 1395
$ws.Range("K136").Value = 53771220
$ws.Range("L136").Value = 4185
$ws.Range("M136").Value = -53768670
$ws.Range("N136").Value = -9285
